$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D:K -> E:L), carrying
# the quarterly financial figures for the latest reported quarter
# (period ending 2018-09-30 / serial 43373) into the new column while
# pushing the previously-existing quarters one column to the right.
$ws.Columns("D").Insert()

# The inserted column has no number formatting of its own yet; copy the
# formats (date format row 7/38/80, integer format elsewhere) from the
# column immediately to its right (old column D, now E) so the new
# column matches its neighbours.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest quarter's values.
$ws.Cells.Item(7, 4).Value = 43373
$ws.Cells.Item(8, 4).Value = 31700
$ws.Cells.Item(9, 4).Value = 9600
$ws.Cells.Item(10, 4).Value = 22100
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(15, 4).Value = 4600
$ws.Cells.Item(17, 4).Value = 39600
$ws.Cells.Item(18, 4).Value = -7900
$ws.Cells.Item(20, 4).Value = -4200
$ws.Cells.Item(21, 4).Value = -11800
$ws.Cells.Item(22, 4).Value = "NA"
$ws.Cells.Item(23, 4).Value = -12100
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = -12100
$ws.Cells.Item(27, 4).Value = -6800
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = 4200
$ws.Cells.Item(33, 4).Value = -6800
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = -6800
$ws.Cells.Item(38, 4).Value = 43373
$ws.Cells.Item(41, 4).Value = 3200
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 3400
$ws.Cells.Item(44, 4).Value = 1100
$ws.Cells.Item(45, 4).Value = 1600
$ws.Cells.Item(46, 4).Value = 9300
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 140200
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 200
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 149700
$ws.Cells.Item(57, 4).Value = 9800
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(59, 4).Value = 16700
$ws.Cells.Item(60, 4).Value = 26600
$ws.Cells.Item(61, 4).Value = 176600
$ws.Cells.Item(62, 4).Value = 51800
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 281800
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = -9900
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = -132100
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43373
$ws.Cells.Item(81, 4).Value = -6800
$ws.Cells.Item(83, 4).Value = "NA"
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = -6300
$ws.Cells.Item(91, 4).Value = -5100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -4800
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = 7600
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 4).Value = -3500
